$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.04358184337616
$ws.Range("B1").Value = 0.9962736368179321
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.067568302154541
$ws.Range("E1").Value = 1.035061359405518
